$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D..K shift right to F..M).
# xlShiftToRight = -4161, CopyOrigin xlFormatFromRightOrBelow = 1
$ws.Columns("D:E").Insert(-4161, 1)

# The Insert() above does not reliably carry the number formatting onto the
# two newly inserted columns, so explicitly copy formatting from column F
# (the former column D) across the new D:E columns. Bound the range to the
# sheet's used rows (5-102) so we don't touch/extend unused rows.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 233200
$ws.Range("E8").Value = 232400
$ws.Range("D9").Value = 119700
$ws.Range("E9").Value = 119100
$ws.Range("D10").Value = 113500
$ws.Range("E10").Value = 113300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 10200
$ws.Range("E14").Value = 7200
$ws.Range("D15").Value = 43400
$ws.Range("E15").Value = 43800
$ws.Range("D17").Value = 192000
$ws.Range("E17").Value = 195800
$ws.Range("D18").Value = 41200
$ws.Range("E18").Value = 36600
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 84500
$ws.Range("E21").Value = 80400
$ws.Range("D22").Value = 23900
$ws.Range("E22").Value = 23500
$ws.Range("D23").Value = 17200
$ws.Range("E23").Value = 13100
$ws.Range("D24").Value = 4200
$ws.Range("E24").Value = 3100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 13000
$ws.Range("E26").Value = 10000
$ws.Range("D27").Value = 12700
$ws.Range("E27").Value = 9700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 12700
$ws.Range("E33").Value = 9700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 12700
$ws.Range("E35").Value = 9700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 5600
$ws.Range("E41").Value = 3400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 163200
$ws.Range("E43").Value = 141800
$ws.Range("D44").Value = 76300
$ws.Range("E44").Value = 77500
$ws.Range("D45").Value = 11000
$ws.Range("E45").Value = 10000
$ws.Range("D46").Value = 256200
$ws.Range("E46").Value = 232700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2171000
$ws.Range("E48").Value = 2165800
$ws.Range("D49").Value = 52400
$ws.Range("E49").Value = 56300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 72900
$ws.Range("E52").Value = 76600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2552500
$ws.Range("E54").Value = 2531400
$ws.Range("D57").Value = 54900
$ws.Range("E57").Value = 71000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 95800
$ws.Range("E59").Value = 88100
$ws.Range("D60").Value = 150700
$ws.Range("E60").Value = 159000
$ws.Range("D61").Value = 1529500
$ws.Range("E61").Value = 1515700
$ws.Range("D62").Value = 30700
$ws.Range("E62").Value = 28900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1710900
$ws.Range("E66").Value = 1703600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -2263700
$ws.Range("E72").Value = -2259500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 841600
$ws.Range("E76").Value = 827800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 12700
$ws.Range("E81").Value = 9700
$ws.Range("D83").Value = 43400
$ws.Range("E83").Value = 43800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 55200
$ws.Range("E89").Value = 65500
$ws.Range("D91").Value = -77900
$ws.Range("E91").Value = -109000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -68100
$ws.Range("E94").Value = -104100
$ws.Range("D96").Value = -17200
$ws.Range("E96").Value = -17100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 15000
$ws.Range("E100").Value = 38500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 2200
$ws.Range("E102").Value = -100


# Data correction: row 94 (Changes In Inventories), the value that lands in
# the new column H (old column F) is corrected from -45200 to -44900.
$ws.Range("H94").Value = -44900
